$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M11").Value = -289.625
$ws.Range("I11").Value = 429.625
$ws.Range("H11").Value = 429.625
$ws.Range("K11").Value = 429.625
$ws.Range("H33").Value = 237.78947
$ws.Range("I33").Value = 223.0625
$ws.Range("M33").Value = 5.9375
$ws.Range("K33").Value = 223.0625
$ws.Range("M42").Value = -212.71429
$ws.Range("K42").Value = 442.71429
$ws.Range("H42").Value = 941.5333000000001
$ws.Range("I42").Value = 147.57143
$ws.Range("L55").Value = 2709.75
$ws.Range("N55").Value = -3137.75
$ws.Range("M55").Value = 114.454544
$ws.Range("I55").Value = 99.545456
$ws.Range("J55").Value = 2709.75
$ws.Range("H55").Value = 501.1154
$ws.Range("K55").Value = 99.545456
$ws.Range("J70").Value = 1389.0834
$ws.Range("L70").Value = 4167.2502
$ws.Range("N70").Value = -4707.2502
$ws.Range("H70").Value = 50610812
$ws.Range("J73").Value = 1389.0834
$ws.Range("L73").Value = 4167.2502
$ws.Range("N73").Value = -6039.2502
$ws.Range("H73").Value = 50610812
$ws.Range("I74").Value = 5015.4
$ws.Range("H74").Value = 5015.4
$ws.Range("M74").Value = -4079.4
$ws.Range("K74").Value = 5015.4
$ws.Range("I77").Value = 5015.4
$ws.Range("M77").Value = -20397
$ws.Range("H77").Value = 5015.4
$ws.Range("K77").Value = 25077
$ws.Range("M82").Value = -17718.5
$ws.Range("I82").Value = 6041.5
$ws.Range("H82").Value = 14079
$ws.Range("K82").Value = 18124.5
$ws.Range("K85").Value = 18124.5
$ws.Range("M85").Value = -16720.5
$ws.Range("I85").Value = 6041.5
$ws.Range("H85").Value = 14079
$ws.Range("N100").Value = -17301.6
$ws.Range("H100").Value = 13012.125
$ws.Range("I100").Value = 7666.3335
$ws.Range("L100").Value = 16219.6
$ws.Range("K100").Value = 7666.3335
$ws.Range("M100").Value = -7125.3335
$ws.Range("J100").Value = 16219.6
$ws.Range("H104").Value = 129.33333
$ws.Range("K104").Value = 237
$ws.Range("M104").Value = 1510
$ws.Range("I104").Value = 79
$ws.Range("N106").ClearContents()
$ws.Range("L106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("H106").Value = 17999.6
$ws.Range("H107").Value = 530.0625
$ws.Range("I107").Value = 470.07144
$ws.Range("K107").Value = 470.07144
$ws.Range("M107").Value = 1449.92856
$ws.Range("M118").Value = -4328
$ws.Range("H118").Value = 1267.8
$ws.Range("N118").Value = -5663
$ws.Range("I118").Value = 1995
$ws.Range("J118").Value = 783
$ws.Range("L118").Value = 2349
$ws.Range("K118").Value = 5985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").Value = -2587.3333
$ws.Range("J2").Value = 2361.3333
$ws.Range("H2").Value = 2687.1853
$ws.Range("L2").Value = 2361.3333
$ws.Range("M2").Value = -2667.2856
$ws.Range("I2").Value = 2780.2856
$ws.Range("K2").Value = 2780.2856
$ws.Range("M32").Value = -2602.3447
$ws.Range("I32").Value = 2889.3447
$ws.Range("H32").Value = 3162.4067
$ws.Range("K32").Value = 2889.3447
$ws.Range("J43").Value = 34535.43
$ws.Range("H43").Value = 35570.89
$ws.Range("N43").Value = -35161.43
$ws.Range("L43").Value = 34535.43
$ws.Range("L61").Value = 3334664
$ws.Range("N61").Value = -3335088
$ws.Range("H61").Value = 31432284
$ws.Range("J61").Value = 3334664
$ws.Range("J74").Value = 7332.6665
$ws.Range("N74").Value = -9080.666499999999
$ws.Range("I74").Value = 1837.6
$ws.Range("H74").Value = 3898.25
$ws.Range("M74").Value = -963.5999999999999
$ws.Range("L74").Value = 7332.6665
$ws.Range("K74").Value = 1837.6
$ws.Range("I77").Value = 1837.6
$ws.Range("N77").Value = -45399.3325
$ws.Range("M77").Value = -4820
$ws.Range("L77").Value = 36663.3325
$ws.Range("H77").Value = 3898.25
$ws.Range("K77").Value = 9188
$ws.Range("J77").Value = 7332.6665
$ws.Range("L88").Value = 2459.8572
$ws.Range("J88").Value = 2459.8572
$ws.Range("I88").Value = 1248
$ws.Range("N88").Value = -3271.8572
$ws.Range("M88").Value = -842
$ws.Range("H88").Value = 2019.1818
$ws.Range("K88").Value = 1248
$ws.Range("H91").Value = 2019.1818
$ws.Range("J91").Value = 2459.8572
$ws.Range("L91").Value = 2459.8572
$ws.Range("I91").Value = 1248
$ws.Range("K91").Value = 1248
$ws.Range("M91").Value = 156
$ws.Range("N91").Value = -5267.8572
$ws.Range("J97").Value = 1320
$ws.Range("N97").Value = -2312
$ws.Range("L97").Value = 1320
$ws.Range("I97").Value = 1026.7826
$ws.Range("M97").Value = -530.7826
$ws.Range("H97").Value = 1060.6154
$ws.Range("K97").Value = 1026.7826
$ws.Range("M102").Value = -214.6666
$ws.Range("K102").Value = 1836.6666
$ws.Range("H102").Value = 2304.2
$ws.Range("I102").Value = 1836.6666
$ws.Range("K110").Value = 6899.5
$ws.Range("I110").Value = 6899.5
$ws.Range("M110").Value = -4854.5
$ws.Range("H110").Value = 6419.6
$ws.Range("H116").Value = 2687.1853
$ws.Range("K116").Value = 2780.2856
$ws.Range("J116").Value = 2361.3333
$ws.Range("I116").Value = 2780.2856
$ws.Range("N116").Value = -6949.3333
$ws.Range("M116").Value = -486.2856000000002
$ws.Range("L116").Value = 2361.3333
$ws.Range("H132").Value = 5563638
$ws.Range("I132").Value = 5419.615
$ws.Range("K132").Value = 16258.845
$ws.Range("M132").Value = -13728.845
$ws.Range("N136").Value = -10009092
$ws.Range("L136").Value = 10003992
$ws.Range("H136").Value = 31432284
$ws.Range("J136").Value = 3334664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 2361.3333
$ws.Range("N3").Value = -2589.3333
$ws.Range("H3").Value = 2687.1853
$ws.Range("I3").Value = 2780.2856
$ws.Range("L3").Value = 2361.3333
$ws.Range("M3").Value = -2666.2856
$ws.Range("K3").Value = 2780.2856
$ws.Range("K94").Value = 2309
$ws.Range("J94").Value = 1862
$ws.Range("I94").Value = 2309
$ws.Range("N94").Value = -2764
$ws.Range("L94").Value = 1862
$ws.Range("M94").Value = -1858
$ws.Range("H94").Value = 2207.4092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M7").Value = -636
$ws.Range("H7").Value = 437
$ws.Range("I7").Value = 749
$ws.Range("K7").Value = 749
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2150
$ws.Range("H22").Value = 2674.875
$ws.Range("I22").Value = 2500
$ws.Range("N22").Value = -3549.75
$ws.Range("J22").Value = 2849.75
$ws.Range("L22").Value = 2849.75
$ws.Range("M31").Value = -45457161
$ws.Range("H31").Value = 27780612
$ws.Range("K31").Value = 45457456
$ws.Range("I31").Value = 45457456
$ws.Range("H34").Value = 27780612
$ws.Range("M34").Value = -45457254
$ws.Range("K34").Value = 45457456
$ws.Range("I34").Value = 45457456
$ws.Range("N99").Value = -8895
$ws.Range("H99").Value = 21023
$ws.Range("L99").Value = 5899
$ws.Range("J99").Value = 5899
$ws.Range("L126").Value = 17697
$ws.Range("H126").Value = 21023
$ws.Range("J126").Value = 5899
$ws.Range("N126").Value = -22637
$ws.Range("H132").Value = 2887.5925
$ws.Range("L132").Value = 10207.875
$ws.Range("J132").Value = 3402.625
$ws.Range("N132").Value = -15267.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L55").Value = 40699.2
$ws.Range("N55").Value = -41053.2
$ws.Range("M55").Value = -6475.6362
$ws.Range("I55").Value = 2217.5454
$ws.Range("J55").Value = 13566.4
$ws.Range("H55").Value = 5764.0625
$ws.Range("K55").Value = 6652.6362
$ws.Range("N92").Value = -2938.5
$ws.Range("H92").Value = 147.5
$ws.Range("J92").Value = 147.5
$ws.Range("L92").Value = 442.5
$ws.Range("M115").Value = -1075
$ws.Range("H115").Value = 17041.5
$ws.Range("K115").Value = 2250
$ws.Range("I115").Value = 750
$ws.Range("L122").Value = 32735.25
$ws.Range("H122").Value = 58192.168
$ws.Range("N122").Value = -37635.25
$ws.Range("J122").Value = 3637.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L80").Value = 4195.25
$ws.Range("J80").Value = 4195.25
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("H80").Value = 4195.25
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6191.25
$ws.Range("J83").Value = 4195.25
$ws.Range("L83").Value = 20976.25
$ws.Range("K83").Value = 0
$ws.Range("N83").Value = -30960.25
$ws.Range("H83").Value = 4195.25
$ws.Range("M83").ClearContents()
$ws.Range("I83").Value = 0
$ws.Range("J97").Value = 164.75
$ws.Range("N97").Value = -1156.75
$ws.Range("L97").Value = 164.75
$ws.Range("I97").Value = 600.55554
$ws.Range("M97").Value = -104.55554
$ws.Range("H97").Value = 521.3182
$ws.Range("K97").Value = 600.55554
$ws.Range("N102").Value = -5834.8572
$ws.Range("M102").Value = -1102.5
$ws.Range("K102").Value = 2724.5
$ws.Range("H102").Value = 2699.2163
$ws.Range("J102").Value = 2590.8572
$ws.Range("L102").Value = 2590.8572
$ws.Range("I102").Value = 2724.5
$ws.Range("H132").Value = 14288785
$ws.Range("L132").Value = 60008400
$ws.Range("I132").Value = 3749.5
$ws.Range("J132").Value = 20002800
$ws.Range("K132").Value = 11248.5
$ws.Range("N132").Value = -60013460
$ws.Range("M132").Value = -8718.5
$ws.Range("N135").Value = -118989
$ws.Range("L135").Value = 108849
$ws.Range("J135").Value = 108849
$ws.Range("H135").Value = 108849

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -4668.857
$ws.Range("J7").Value = 8249.5
$ws.Range("H7").Value = 5551.6665
$ws.Range("L7").Value = 8249.5
$ws.Range("I7").Value = 4780.857
$ws.Range("K7").Value = 4780.857
$ws.Range("N7").Value = -8473.5
$ws.Range("H29").Value = 75000
$ws.Range("M29").ClearContents()
$ws.Range("K29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("H46").Value = 2749.6667
$ws.Range("L46").Value = 3416.6667
$ws.Range("N46").Value = -3792.6667
$ws.Range("J46").Value = 3416.6667
$ws.Range("H93").Value = 3273285.2
$ws.Range("M93").Value = -1325.7273
$ws.Range("N93").Value = -9272086
$ws.Range("K93").Value = 2573.7273
$ws.Range("L93").Value = 9269590
$ws.Range("I93").Value = 2573.7273
$ws.Range("J93").Value = 9269590
$ws.Range("I122").Value = 3087.5366
$ws.Range("H122").Value = 3235.2444
$ws.Range("M122").Value = -6812.6098
$ws.Range("K122").Value = 9262.6098
$ws.Range("K126").Value = 14342.571
$ws.Range("L126").Value = 24748.5
$ws.Range("M126").Value = -11872.571
$ws.Range("I126").Value = 4780.857
$ws.Range("H126").Value = 5551.6665
$ws.Range("J126").Value = 8249.5
$ws.Range("N126").Value = -29688.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L62").Value = 15273
$ws.Range("K62").Value = 6947
$ws.Range("I62").Value = 6947
$ws.Range("M62").Value = -6323
$ws.Range("N62").Value = -16521
$ws.Range("J62").Value = 15273
$ws.Range("H62").Value = 11430.23
$ws.Range("I65").Value = 6947
$ws.Range("L65").Value = 76365
$ws.Range("J65").Value = 15273
$ws.Range("M65").Value = -31615
$ws.Range("H65").Value = 11430.23
$ws.Range("N65").Value = -82605
$ws.Range("K65").Value = 34735
$ws.Range("H81").Value = 1891.2222
$ws.Range("J81").Value = 3348
$ws.Range("N81").Value = -8818
$ws.Range("K81").Value = 3418.25
$ws.Range("L81").Value = 6696
$ws.Range("I81").Value = 1709.125
$ws.Range("M81").Value = -2357.25
$ws.Range("J84").Value = 3348
$ws.Range("I84").Value = 1709.125
$ws.Range("M84").Value = -11787.25
$ws.Range("N84").Value = -44088
$ws.Range("L84").Value = 33480
$ws.Range("K84").Value = 17091.25
$ws.Range("H84").Value = 1891.2222
$ws.Range("H107").Value = 4504.8667
$ws.Range("J107").Value = 5169.2856
$ws.Range("L107").Value = 15507.8568
$ws.Range("I107").Value = 2954.5557
$ws.Range("N107").Value = -19347.8568
$ws.Range("K107").Value = 8863.667099999999
$ws.Range("M107").Value = -6943.667099999999
$ws.Range("I122").Value = 2667.9167
$ws.Range("H122").Value = 2567.6667
$ws.Range("M122").Value = -5553.750100000001
$ws.Range("K122").Value = 8003.750100000001
